$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, pushing existing rows 39:121 down to 40:122
$ws.Rows("39:39").Insert()

# Populate the newly inserted row 39 with the new data record
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 45100
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 100112010
$ws.Range("G39").Value = "Achicoria"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 65
$ws.Range("K39").Value = 8000
$ws.Range("L39").Value = 8000
$ws.Range("M39").Value = 8000
$ws.Range("N39").Value = "$/caja 18 unidades"
$ws.Range("O39").Value = "Región Metropolitana"
$ws.Range("P39").Value = 444
$ws.Range("Q39").Value = 18
$ws.Range("R39").Value = "Hortaliza"
